$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title text updates (Volume/Number and Report Week date range) ---
$ws.Range("A8").Value = "Volume 32   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/13/2025  Through  10/19/2025"

# --- Model cells used as formatting sources for cells whose type changes ---
# D14 = General/text style (s=13); D16 = integer style (s=14); E16 = one-decimal style (s=15)
$styleText = $ws.Range("D14")
$styleInt  = $ws.Range("D16")
$styleDec  = $ws.Range("E16")

# --- Apply formatting (style) changes first, before setting any values ---
$styleInt.Copy()
$ws.Range("C15").PasteSpecial(-4122)
$styleInt.Copy()
$ws.Range("F15").PasteSpecial(-4122)
$styleInt.Copy()
$ws.Range("D18").PasteSpecial(-4122)
$styleDec.Copy()
$ws.Range("E18").PasteSpecial(-4122)
$styleText.Copy()
$ws.Range("C20").PasteSpecial(-4122)
$styleInt.Copy()
$ws.Range("C22").PasteSpecial(-4122)
$styleInt.Copy()
$ws.Range("D23").PasteSpecial(-4122)
$styleDec.Copy()
$ws.Range("E23").PasteSpecial(-4122)
$styleText.Copy()
$ws.Range("D25").PasteSpecial(-4122)
$styleText.Copy()
$ws.Range("E25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Apply all cell value changes ---
$ws.Range("N14").Value = -80
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 28
$ws.Range("K15").Value = 40
$ws.Range("L15").Value = 115.384615384615
$ws.Range("M15").Value = 115.384615384615
$ws.Range("N15").Value = -20
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 79
$ws.Range("J16").Value = 111
$ws.Range("K16").Value = -28.828828828828
$ws.Range("L16").Value = -7.058823529411
$ws.Range("M16").Value = -20.202020202020
$ws.Range("N16").Value = -85.531135531135
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -27.272727272727
$ws.Range("F17").Value = 32
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = 14.285714285714
$ws.Range("I17").Value = 311
$ws.Range("J17").Value = 288
$ws.Range("K17").Value = 7.986111111111
$ws.Range("L17").Value = 43.317972350230
$ws.Range("M17").Value = 142.96875
$ws.Range("N17").Value = -25.060240963855
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 57
$ws.Range("J18").Value = 58
$ws.Range("K18").Value = -1.724137931034
$ws.Range("L18").Value = 1.785714285714
$ws.Range("M18").Value = -47.706422018348
$ws.Range("N18").Value = -91.190108191653
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 100
$ws.Range("G19").Value = 18
$ws.Range("H19").Value = -11.111111111111
$ws.Range("I19").Value = 142
$ws.Range("J19").Value = 140
$ws.Range("K19").Value = 1.428571428571
$ws.Range("L19").Value = -8.974358974358
$ws.Range("M19").Value = 79.746835443038
$ws.Range("N19").Value = -45.593869731800
$ws.Range("C20").Value = "0"
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -100
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -62.5
$ws.Range("J20").Value = 63
$ws.Range("K20").Value = -52.380952380952
$ws.Range("M20").Value = -55.882352941176
$ws.Range("N20").Value = -91.957104557640
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -19.047619047619
$ws.Range("F21").Value = 67
$ws.Range("G21").Value = 75
$ws.Range("H21").Value = -10.666666666666
$ws.Range("I21").Value = 650
$ws.Range("J21").Value = 683
$ws.Range("K21").Value = -4.831625183016
$ws.Range("L21").Value = 12.068965517241
$ws.Range("M21").Value = 29.740518962075
$ws.Range("N21").Value = -71.640488656195
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 9
$ws.Range("K22").Value = -10
$ws.Range("L22").Value = -35.714285714285
$ws.Range("M22").Value = 12.5
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -100
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = -28.571428571428
$ws.Range("J23").Value = 70
$ws.Range("K23").Value = -10
$ws.Range("L23").Value = -5.970149253731
$ws.Range("M23").Value = 85.294117647058
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 58.333333333333
$ws.Range("F24").Value = 98
$ws.Range("G24").Value = 62
$ws.Range("H24").Value = 58.064516129032
$ws.Range("I24").Value = 702
$ws.Range("J24").Value = 538
$ws.Range("K24").Value = 30.483271375464
$ws.Range("L24").Value = 51.293103448275
$ws.Range("M24").Value = 138.775510204082
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = "0"
$ws.Range("E25").Value = "***.*"
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 1250
$ws.Range("I25").Value = 130
$ws.Range("K25").Value = 242.105263157895
$ws.Range("L25").Value = 60.493827160493
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = 22.222222222222
$ws.Range("F26").Value = 55
$ws.Range("G26").Value = 41
$ws.Range("H26").Value = 34.146341463414
$ws.Range("I26").Value = 409
$ws.Range("J26").Value = 352
$ws.Range("K26").Value = 16.193181818181
$ws.Range("L26").Value = 18.550724637681
$ws.Range("M26").Value = 22.089552238806
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 34
$ws.Range("K27").Value = 9.677419354838
$ws.Range("L27").Value = 36
$ws.Range("C28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -20
$ws.Range("I28").Value = 33
$ws.Range("J28").Value = 37
$ws.Range("K28").Value = -10.810810810810
$ws.Range("L28").Value = -2.941176470588
$ws.Range("G29").Value = 2
$ws.Range("G30").Value = 2
$ws.Range("F31").Value = 4
$ws.Range("I31").Value = 6
$ws.Range("K31").Value = 20
$ws.Range("L31").Value = 200
